$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.315621137619019
$ws.Range("B1").Value = 1.625413060188293
$ws.Range("C1").Value = 2.259082078933716
$ws.Range("D1").Value = 4.833589553833008
$ws.Range("E1").Value = 4.305968761444092
